$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at row 272 (shifting existing rows, e.g. IAD, down by one)
$ws.Rows.Item(272).Insert()

# Copy the formatting from the row that just shifted down (originally row 272,
# now row 273, e.g. the IAD row) onto the newly inserted row so the "colo"
# cell keeps its bold/bordered/centered style.
$ws.Range("A273").Copy()
$ws.Range("A272").PasteSpecial(-4122)

# Populate the new row with the Chengdu, China colo entry.
$ws.Range("A272").Value = "CTU"
$ws.Range("B272").Value = "Chengdu, China"
$ws.Range("C272").Value = "Asia"
$ws.Range("D272").Value = "Chengdu"
$ws.Range("E272").Value = "China"
$ws.Range("F272").Value = "CN"

# Latitude/longitude are left blank for this new entry.
$ws.Range("G272").ClearContents()
$ws.Range("H272").ClearContents()
